$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.419.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -2.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4295"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3701"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8665"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.697"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.359"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07068"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008905"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.429.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.168"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.062.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.000"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.151"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.291"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08866"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7663"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.483"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.909"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01961"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05291"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.179"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.872"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5083"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1674"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.616"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4742"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06424"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.666"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.832"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
